$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    "
$find.Replacement.ClearFormatting()
$find.Replacement.Text = ""
$find.Execute(
    $find.Text, $true, $false, $false, $false, $false,
    $true, 1, $false, $find.Replacement.Text, 2
)
